$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V3").Value = 1133000
$ws.Range("V5").Value = 1013000
$ws.Range("V6").Value = 1770000
$ws.Range("D9").Value = 805900
$ws.Range("V9").Value = 805900
$ws.Range("V10").Value = 574000
$ws.Range("E14").Value = "https://carsklad-174.ru/auto/changan/cs55/i/suv-5d"
$ws.Range("V14").Value = 759900
$ws.Range("V15").Value = 1233900
$ws.Range("V17").Value = 1084300
$ws.Range("V18").Value = 2079900
$ws.Range("V21").Value = 1319900
$ws.Range("E24").Value = "https://carsklad-174.ru/auto/changan/uni-k/i/suv-5d"
$ws.Range("V24").Value = 2039900
$ws.Range("V27").Value = 1623900
$ws.Range("V30").Value = 560900
$ws.Range("E35").Value = "https://carsklad-174.ru/auto/chery/tiggo-7-pro/i/suv-5d"
$ws.Range("V35").Value = 1087900
$ws.Range("V36").Value = 1299000
$ws.Range("E40").Value = "https://carsklad-174.ru/auto/chery/tiggo-8-pro/i/suv-5d"
$ws.Range("V40").Value = 1179900
$ws.Range("D42").Value = 1549900
$ws.Range("V42").Value = 1549900
$ws.Range("V45").Value = 537900
$ws.Range("V46").Value = 487900
$ws.Range("V47").Value = 537900
$ws.Range("D52").Value = 754000
$ws.Range("E52").Value = "https://carsklad-174.ru/auto/dfm/580/ii/suv-5d"
$ws.Range("V52").Value = 754000
$ws.Range("E56").Value = "https://carsklad-174.ru/auto/exeed/txl/i/suv-5d"
$ws.Range("V56").Value = 1399900
$ws.Range("V64").Value = 614300
$ws.Range("V65").Value = 756000
$ws.Range("V71").Value = 864900
$ws.Range("V72").Value = 1529000
$ws.Range("V73").Value = 1082596
$ws.Range("E75").Value = "https://carsklad-174.ru/auto/geely/coolray/i/suv-5d"
$ws.Range("V75").Value = 619990
$ws.Range("V76").Value = 1219990
$ws.Range("V81").Value = 2689900
$ws.Range("E84").Value = "https://carsklad-174.ru/auto/geely/tugella/i-restyling/suv-5d"
$ws.Range("V84").Value = 1859990
$ws.Range("V89").Value = 1674000
$ws.Range("E90").Value = "https://carsklad-174.ru/auto/haval/dargo/i/x"
$ws.Range("V90").Value = 1959000
$ws.Range("E91").Value = "https://carsklad-174.ru/auto/haval/f_7/i/suv-5d"
$ws.Range("V91").Value = 1119000
$ws.Range("E93").Value = "https://carsklad-174.ru/auto/haval/f_7_x/i/suv-5d"
$ws.Range("V93").Value = 1199000
$ws.Range("E99").Value = "https://carsklad-174.ru/auto/haval/h9/i/suv-5d"
$ws.Range("V99").Value = 1803000
$ws.Range("E101").Value = "https://carsklad-174.ru/auto/haval/jolion/i/suv-5d"
$ws.Range("V101").Value = 689900
$ws.Range("V102").Value = 1199000
$ws.Range("V104").Value = 671000
$ws.Range("V106").Value = 863000
$ws.Range("V110").Value = 460000
$ws.Range("E111").Value = "https://carsklad-174.ru/auto/hyundai/sonata/viii/sedan"
$ws.Range("V111").Value = 1432000
$ws.Range("E113").Value = "https://carsklad-174.ru/auto/hyundai/tucson/iv/suv-5d"
$ws.Range("V113").Value = 1429000
$ws.Range("E116").Value = "https://carsklad-174.ru/auto/jac/j7/i/liftback"
$ws.Range("V116").Value = 709000
$ws.Range("E133").Value = "https://carsklad-174.ru/auto/jetta/va3/i/sedan"
$ws.Range("V133").Value = 840000
$ws.Range("E134").Value = "https://carsklad-174.ru/auto/jetta/vs5/i/suv-5d"
$ws.Range("V134").Value = 1140000
$ws.Range("V136").Value = 1053000
$ws.Range("V137").Value = 1317900
$ws.Range("V138").Value = 1350000
$ws.Range("V141").Value = 859900
$ws.Range("V142").Value = 969900
$ws.Range("V143").Value = 943990
$ws.Range("V144").Value = 1065300
$ws.Range("E147").Value = "https://carsklad-174.ru/auto/kia/picanto/iii-restyling/hatchback-5d"
$ws.Range("V147").Value = 619900
$ws.Range("V149").Value = 460900
$ws.Range("V150").Value = 590900
$ws.Range("V152").Value = 920300
$ws.Range("V154").Value = 571200
$ws.Range("E158").Value = "https://carsklad-174.ru/auto/lada/granta/i-restyling/cross"
$ws.Range("V158").Value = 304400
$ws.Range("V159").Value = 263500
$ws.Range("V160").Value = 267300
$ws.Range("V161").Value = 252900
$ws.Range("E162").Value = "https://carsklad-174.ru/auto/lada/granta/i-restyling/drive-active"
$ws.Range("V162").Value = 412500
$ws.Range("V167").Value = 273900
$ws.Range("V168").Value = 343900
$ws.Range("V169").Value = 379300
$ws.Range("V179").Value = 355200
$ws.Range("V180").Value = 439500
$ws.Range("V182").Value = 428000
$ws.Range("V183").Value = 585900
$ws.Range("V187").Value = 375900
$ws.Range("V188").Value = 460900
$ws.Range("V192").Value = 380900
$ws.Range("V193").Value = 379300
$ws.Range("V207").Value = 2434000
$ws.Range("E209").Value = "https://carsklad-174.ru/auto/nissan/qashqai/2-rest/suv-5d"
$ws.Range("V209").Value = 880000
$ws.Range("E210").Value = "https://carsklad-174.ru/auto/nissan/terrano/iii/suv-5d"
$ws.Range("V210").Value = 619000
$ws.Range("V211").Value = 1129000
$ws.Range("V212").Value = 1263900
$ws.Range("V217").Value = 853000
$ws.Range("V218").Value = 470000
$ws.Range("E220").Value = "https://carsklad-174.ru/auto/renault/kaptur/i-rest/suv-5d"
$ws.Range("V220").Value = 614000
$ws.Range("V222").Value = 390000
$ws.Range("V223").Value = 498990
$ws.Range("V225").Value = 390000
$ws.Range("V226").Value = 394000
$ws.Range("V228").Value = 949000
$ws.Range("V231").Value = 3058000
$ws.Range("V234").Value = 2745000
$ws.Range("V235").Value = 839800
$ws.Range("V238").Value = 420000
$ws.Range("V239").Value = 1206000
$ws.Range("V249").Value = 3893000
$ws.Range("V260").Value = 514300
$ws.Range("V261").Value = 1222900
$ws.Range("V263").Value = 1283900
$ws.Range("E72").Value = "https://carsklad-174.ru/auto/geely/atlas/ii/suv-5d"
$ws.Range("E76").Value = "https://carsklad-174.ru/auto/geely/coolray/i-restyling/suv-5d"
$ws.Range("E102").Value = "https://carsklad-174.ru/auto/haval/jolion/i-rest/suv-5d"
$ws.Range("E231").Value = "https://carsklad-174.ru/auto/skoda/kodiaq/i-rest/laurin"
